$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 13 ("Programa resumido:" / duplicate "5816812 - João Paulo Alves Silva"
# text) is removed entirely. Deleting it shifts every following row up by one and
# keeps the custom row heights of the old rows 14-24 intact (they become rows 13-23),
# which matches the target layout automatically.
$ws.Rows.Item(13).Delete()

# --- Row 10: the "Objetivos:" answer becomes the professor's name ---
$ws.Range("B10").Value = "5816812 - João Paulo Alves Silva"
$ws.Range("C10").Value = "5816812 - João Paulo Alves Silva"

# --- Row 13 (was old row 14): "Programa resumido:" / "Semestral" ---
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# --- Row 14 (was old row 15): "Short syllabus:" only ---
$ws.Range("A14").Value = "Short syllabus:"

# --- Row 15 (was old row 16): "Programa:" / "01/01/2012" ---
# NOTE: "01/01/2012" looks like a date, and a plain .Value assignment would get
# auto-converted into a date serial number by Excel. To keep it as plain text
# (matching the shared-string cell type used elsewhere in the sheet) without
# disturbing the cell's existing style (which would happen if we touched
# NumberFormat directly), build the text via a formula on a scratch cell, copy
# it, and paste just the computed value into the target cells.
$ws.Range("A15").Value = "Programa:"
$ws.Range("Z1").Formula = "=""01/01/2012"""
$ws.Range("Z1").Copy()
$ws.Range("B15").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("C15").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("Z1").EntireColumn.Delete()

# --- Row 16 (was old row 17): "Syllabus:" only ---
$ws.Range("A16").Value = "Syllabus:"

# --- Row 17 (was old row 18): "Avaliação:" only ---
$ws.Range("A17").Value = "Avaliação:"

# --- Row 18 (was old row 19): "Método:" / "5816812 - João Paulo Alves Silva" ---
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "5816812 - João Paulo Alves Silva"
$ws.Range("C18").Value = "5816812 - João Paulo Alves Silva"

# --- Row 19 (was old row 20): "Critério:" / "O desenvolvimento..." ---
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "O desenvolvimento da disciplina será baseado em leituras, aula expositiva, discussão e resolução de estudos de caso e resolução de exercícios."
$ws.Range("C19").Value = "O desenvolvimento da disciplina será baseado em leituras, aula expositiva, discussão e resolução de estudos de caso e resolução de exercícios."

# --- Row 20 (was old row 21): "Norma de recuperação:" / "Provas e trabalhos." ---
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "Provas e trabalhos."
$ws.Range("C20").Value = "Provas e trabalhos."

# --- Row 21 (was old row 22): "Bibliografia:" / "Prova única com nota maior ou igual a 5,0 (cinco)." ---
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "Prova única com nota maior ou igual a 5,0 (cinco)."
$ws.Range("C21").Value = "Prova única com nota maior ou igual a 5,0 (cinco)."

# --- Row 22 (was old row 23): "Requisitos:" only ---
$ws.Range("A22").Value = "Requisitos:"

# --- Row 23 (was old row 24): requirement text, no A cell ---
$ws.Range("B23").Value = "LOQ4064 -  Engenharia de Processos Quimicos I  (Requisito fraco)`n"
$ws.Range("C23").Value = "LOQ4064 -  Engenharia de Processos Quimicos I  (Requisito fraco)`n"
